$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant cell values (text/number changes driven by the new simulation state)
$ws.Range("L2").Value = "x"
$ws.Range("D5").Value = "Scott"
$ws.Range("J5").Value = "x"
$ws.Range("D6").Value = "Sally"
$ws.Range("G6").Value = "Sally"
$ws.Range("A7").Value = "x"
$ws.Range("F7").Value = "x"
$ws.Range("O7").Value = "x"
$ws.Range("I8").Value = "x"
$ws.Range("F9").Value = "Tomas"
$ws.Range("M10").Value = "Jim"
$ws.Range("H11").Value = "Lillian"
$ws.Range("F13").Value = "x"
$ws.Range("G14").Value = "Joe"
$ws.Range("H14").Value = "Diane"
$ws.Range("B15").Value = "x"
$ws.Range("I15").Value = "Scott"
$ws.Range("M15").Value = "x"
$ws.Range("M16").Value = "x"
$ws.Range("D18").Value = "Diane"
$ws.Range("F18").Value = "Tomas"
$ws.Range("G18").Value = "Scott"
$ws.Range("B19").Value = "6,13"
$ws.Range("C19").Value = "12,9"
$ws.Range("D19").Value = "5,8"
$ws.Range("E19").Value = "6,5"
$ws.Range("F19").Value = "7,13"
$ws.Range("G19").Value = "3,5"
$ws.Range("H19").Value = "8,14"
$ws.Range("I19").Value = "7,10"
$ws.Range("J19").Value = "3,4"
$ws.Range("B20").Value = "resting"
$ws.Range("C20").Value = "potential_client"
$ws.Range("D20").Value = "client"
$ws.Range("E20").Value = "resting"
$ws.Range("F20").Value = "host"
$ws.Range("G20").Value = "host"
$ws.Range("H20").Value = "potential_client"
$ws.Range("I20").Value = "resting"
$ws.Range("C21").Value = "Sprint"
$ws.Range("D21").Value = "ATnT"
$ws.Range("F21").Value = "ATnT"
$ws.Range("G21").Value = "Verizon"
$ws.Range("B22").Value = 2
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 4
$ws.Range("C24").Value = "Joe"
$ws.Range("D24").Value = "Joe"
$ws.Range("E24").Value = "Jim"
$ws.Range("F24").Value = "Joe"
$ws.Range("G24").Value = "Tomas"
$ws.Range("H24").Value = "Joe"
$ws.Range("I24").Value = "Joe"
$ws.Range("J24").Value = "Tomas"
$ws.Range("B25").Value = "Jim"
$ws.Range("C25").Value = "Jim"
$ws.Range("E25").Value = "Tomas"
$ws.Range("F25").Value = "Jim"
$ws.Range("H25").Value = "Jim"
$ws.Range("I25").Value = "Jim"
$ws.Range("J25").Value = "Sally"
$ws.Range("B26").Value = "Tomas"
$ws.Range("C26").Value = "Sally"
$ws.Range("D26").Value = "Sally"
$ws.Range("E26").Value = "Sally"
$ws.Range("F26").Value = "Tomas"
$ws.Range("G26").Value = "Sally"
$ws.Range("H26").Value = "Tomas"
$ws.Range("I26").Value = "Tomas"
$ws.Range("B27").Value = "Diane"
$ws.Range("C27").Value = "Diane"
$ws.Range("D27").Value = "Diane"
$ws.Range("E27").Value = "Sally"
$ws.Range("F27").Value = "Diane"
$ws.Range("G27").Value = "Lillian"
$ws.Range("H27").Value = "Diane"
$ws.Range("I27").Value = "Sally"
$ws.Range("J27").Value = "Lillian"
$ws.Range("B28").Value = "Scott"
$ws.Range("D28").Value = "Sally"
$ws.Range("E28").Value = "Lillian"
$ws.Range("G28").Value = "Scott"
$ws.Range("H28").Value = "Scott"
$ws.Range("I28").Value = "Diane"
$ws.Range("J28").Value = "Scott"
$ws.Range("B29").Value = "Lillian"
$ws.Range("C29").Value = "Lillian"
$ws.Range("D29").Value = "Scott"
$ws.Range("E29").Value = "Scott"
$ws.Range("F29").Value = "Lillian"
$ws.Range("H29").Value = "Lillian"
$ws.Range("I29").Value = "Sally"
$ws.Range("D30").Value = "Lillian"
$ws.Range("I30").Value = "Scott"
$ws.Range("D31").Value = "Scott"
$ws.Range("I31").Value = "Lillian"
$ws.Range("I32").Value = "Scott"

# Clear cells that no longer hold data (10th agent column removed, a couple of stray cells cleared)
$ws.Range("K17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("I18").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("K21").ClearContents()
$ws.Range("K22").ClearContents()
$ws.Range("K23").ClearContents()
$ws.Range("K24").ClearContents()
$ws.Range("K25").ClearContents()
$ws.Range("K26").ClearContents()
$ws.Range("J29").ClearContents()

Write-Output "Applied SSC simulation data update"
